# Add a new submission row to the "Tabelle2" table on the active sheet.
# New data: Date=2022-11-26 (serial 44891), Name="221126_cat_test",
# Name Ramp="NewPhoneWhoDis", Hand in="TRUE", By="Maria".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the Excel Table by one row - this also extends ref/autoFilter and
# the sheet dimension once the new row gets content.
$lo = $ws.ListObjects.Item("Tabelle2")
$newRow = $lo.ListRows.Add()

$lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
$newRowIndex = $lastRow
$prevRowIndex = $lastRow - 1

# Date (column A) - same serial/date style as the row above it.
$ws.Cells.Item($newRowIndex, 1).Value = 44891
$ws.Cells.Item($prevRowIndex, 1).Copy()
$ws.Cells.Item($newRowIndex, 1).PasteSpecial(-4122)

# Name Ramp (column C) first, so the new shared string "NewPhoneWhoDis"
# is registered before "221126_cat_test".
$ws.Cells.Item($newRowIndex, 3).Value = "NewPhoneWhoDis"

# Name (column B).
$ws.Cells.Item($newRowIndex, 2).Value = "221126_cat_test"

# Hand in (column D) - copy from the row above so it stays the text
# "TRUE" (shared string) rather than becoming a native boolean.
$ws.Cells.Item($prevRowIndex, 4).Copy($ws.Cells.Item($newRowIndex, 4))

# By (column E).
$ws.Cells.Item($newRowIndex, 5).Value = "Maria"

# Move the selection to the next empty row, matching where Excel leaves
# the cursor after finishing a row entry.
$ws.Range("B25").Select()
